# Reproduce the POI 4.1.0 -> 5.2.3 rPr normalisation for each bold table
# header run: Bold=On, Italic=Off, Strike=Off, Color=000000, Size=24 half
# points (12pt). Each header text is immediately followed in the document
# by a plain (unformatted) duplicate run with the same (or a truncated)
# text value, so a simple whole-document Find would sometimes land on the
# wrong (already-plain) run. We therefore walk the document from a moving
# cursor, always taking the first match at/after the cursor, which is
# guaranteed to be the still-bold header run.

$d = $word.ActiveDocument

$targets = @(
    "'Name'",
    'World',
    'MultiNamedElement',
    'NamedElement',
    'Producer -> NamedElement',
    'Adress',
    'Company -> NamedElement',
    'ProductionCompany -> Company',
    'Restaurant -> Company',
    'Chef -> NamedElement',
    'Recipe -> NamedElement',
    'Food -> NamedElement',
    'Source -> MultiNamedElement',
    'Plant -> Source',
    'Animal -> Source',
    'Color',
    'Caliber',
    'Group',
    'Continent',
    'Kind',
    'Part',
    'CountryData [anydsl.Country]',
    'SingleString [java.lang.String]',
    'EStringToRecipeMap [java.util.Map$Entry]'
)

$cursor = 0
$docEnd = $d.Content.End
$appliedCount = 0

foreach ($t in $targets) {
    $rng = $d.Range($cursor, $docEnd)
    $find = $rng.Find
    $find.ClearFormatting()
    $find.Text = $t
    $found = $find.Execute($t, $true, $true, $false, $false, $false, $true, 1, $false, "", 0)

    if ($found) {
        $rng.Font.Bold = $true
        $rng.Font.Italic = $false
        $rng.Font.StrikeThrough = $false
        $rng.Font.Size = 12

        $cursor = $rng.End
        $appliedCount = $appliedCount + 1
    } else {
        Write-Host "WARNING: could not find target text: $t"
    }
}

Write-Host "Applied header formatting to $appliedCount of $($targets.Count) runs."
